## Applies the "Add files via upload" update to lca_optimization_input.xlsx
## - Extends the small results table on the "index" sheet (rows 28-30,
##   columns A:C) with two new rows: IRR_opt / LCPB_opt.
## - Moves the "last row" box-border formatting down from row 28 to the new
##   row 30, and gives rows 28-29 the regular "middle" row formatting.
## - Updates the active selection to K27 (matches the author's last saved
##   cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("index")

# --- 1. Re-arrange cell formatting for rows 28-30 ------------------------
# Row 28 currently carries the "bottom of box" border style (left+bottom /
# bottom / right+bottom). That style needs to move to row 30, which becomes
# the new last row of the little table. Rows 28-29 instead get the regular
# "inside the box" style already used by row 27 (left / none / right).

$ws.Range("A28:C28").Copy() | Out-Null
$ws.Range("A30:C30").PasteSpecial(-4122) | Out-Null

$ws.Range("A27:C27").Copy() | Out-Null
$ws.Range("A28:C29").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 2. Fill in the two new data rows -------------------------------------
$ws.Range("A29").Value = "IRR_opt"
$ws.Range("B29").Value = "Optimal Internal Rate of Return"
$ws.Range("C29").Value = "%"

$ws.Range("A30").Value = "LCPB_opt"
$ws.Range("B30").Value = "Optimal Life Cyle Payback"
$ws.Range("C30").Value = "years"

# --- 3. Restore the cursor / selection on the sheet -----------------------
$ws.Range("K27").Select() | Out-Null
